# Refresh the "cryptos" price/volume snapshot (Price column D, Volume(1h)
# column E) for rows 2-51. Column D values are leading-apostrophe'd so
# Excel stores them as literal text (matching the sheet's existing
# inline-string cells) instead of auto-coercing numeric-looking strings
# like "0.2930" / "1.000" into floating point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''29.130.06'
$ws.Range("E2").Value = '  +0.03%  '
$ws.Range("D3").Value = '''1.831.15'
$ws.Range("E3").Value = '  -0.32%  '
$ws.Range("D4").Value = '''0.9995'
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '''241.05'
$ws.Range("E5").Value = '  +0.54%  '
$ws.Range("D6").Value = '''0.6624'
$ws.Range("E6").Value = '  -2.54%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").Value = '''0.07388'
$ws.Range("E8").Value = '  -0.59%  '
$ws.Range("D9").Value = '''0.2930'
$ws.Range("D10").Value = '''22.66'
$ws.Range("E10").Value = '  -1.99%  '
$ws.Range("D11").Value = '''0.07730'
$ws.Range("E11").Value = '  +1.16%  '
$ws.Range("D12").Value = '''1.849.88'
$ws.Range("E12").Value = '  +0.57%  '
$ws.Range("D13").Value = '''4.976'
$ws.Range("E13").Value = '  -0.94%  '
$ws.Range("D14").Value = '''0.6661'
$ws.Range("E14").Value = '  -1.94%  '
$ws.Range("D15").Value = '''82.60'
$ws.Range("E15").Value = '  -4.80%  '
$ws.Range("D16").Value = '''6.078'
$ws.Range("E16").Value = '  -1.13%  '
$ws.Range("D17").Value = '''0.000008324'
$ws.Range("E17").Value = '  +1.48%  '
$ws.Range("D18").Value = '''29.146.77'
$ws.Range("E18").Value = '  +0.07%  '
$ws.Range("D19").Value = '''2.088.57'
$ws.Range("E19").Value = '  +0.28%  '
$ws.Range("D20").Value = '''226.48'
$ws.Range("E20").Value = '  -1.21%  '
$ws.Range("D21").Value = '''12.44'
$ws.Range("E21").Value = '  -0.35%  '
$ws.Range("E22").Value = '  +0.14%  '
$ws.Range("D23").Value = '''7.138'
$ws.Range("E23").Value = '  -2.66%  '
$ws.Range("D24").Value = '''1.000'
$ws.Range("E24").Value = '  +0.00%  '
$ws.Range("D25").Value = '''159.41'
$ws.Range("E25").Value = '  -1.23%  '
$ws.Range("D26").Value = '''8.584'
$ws.Range("E26").Value = '  -1.32%  '
$ws.Range("D27").Value = '''0.1397'
$ws.Range("E27").Value = '  -2.20%  '
$ws.Range("D28").Value = '''17.90'
$ws.Range("E28").Value = '  -0.65%  '
$ws.Range("D29").Value = '''1.508'
$ws.Range("E29").Value = '  +0.60%  '
$ws.Range("D30").Value = '''4.098'
$ws.Range("E30").Value = '  -3.43%  '
$ws.Range("D31").Value = '''4.025'
$ws.Range("E31").Value = '  -2.76%  '
$ws.Range("D32").Value = '''1.191'
$ws.Range("E32").Value = '  -0.09%  '
$ws.Range("E33").Value = '  -0.79%  '
$ws.Range("D34").Value = '''1.865'
$ws.Range("E34").Value = '  +1.31%  '
$ws.Range("D35").Value = '''0.7481'
$ws.Range("E35").Value = '  -0.62%  '
$ws.Range("D36").Value = '''1.133'
$ws.Range("E36").Value = '  +0.39%  '
$ws.Range("D37").Value = '''2.642'
$ws.Range("E37").Value = '  -1.56%  '
$ws.Range("D38").Value = '''1.273.87'
$ws.Range("E38").Value = '  -2.97%  '
$ws.Range("E39").Value = '  -1.38%  '
$ws.Range("E40").Value = '  +0.37%  '
$ws.Range("D41").Value = '''0.9288'
$ws.Range("E41").Value = '  -0.68%  '
$ws.Range("D42").Value = '''5.896'
$ws.Range("E42").Value = '  -2.82%  '
$ws.Range("D43").Value = '''0.08462'
$ws.Range("E43").Value = '  +2.60%  '
$ws.Range("E44").Value = '  +0.12%  '
$ws.Range("D45").Value = '''101.74'
$ws.Range("E45").Value = '  -3.01%  '
$ws.Range("D46").Value = '''1.984.85'
$ws.Range("E46").Value = '  -0.07%  '
$ws.Range("E47").Value = '  -0.61%  '
$ws.Range("D48").Value = '''1.757'
$ws.Range("E48").Value = '  -0.59%  '
$ws.Range("E49").Value = '  -1.34%  '
$ws.Range("D50").Value = '''62.93'
$ws.Range("E50").Value = '  -1.75%  '
$ws.Range("D51").Value = '''0.05878'
$ws.Range("E51").Value = '  -0.93%  '
